# testing different plot layouts for #35
# Extend the weekly hours log with two more entries (rows 56 and 57).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (55) down onto the two
# new rows so the date cells keep the same date number format / style.
$ws.Range("A55").Copy()
$ws.Range("A56:A57").PasteSpecial(-4122)

# New date entries (serial dates 45419 = 2024-05-07, 45421 = 2024-05-09)
$ws.Range("A56").Value = 45419
$ws.Range("A57").Value = 45421

# New hours entries
$ws.Range("B56").Value = 0.5
$ws.Range("B57").Value = 1

# Running-total formulas, continuing the existing C-column pattern
$ws.Range("C56").Formula = "=C55+B56"
$ws.Range("C57").Formula = "=C56+B57"

# Move the active selection to the new last cell, matching the saved
# workbook state (B57 selected).
[void]$ws.Range("B57").Select()
